# TOPIC_TABLE COLUMNS TOPIC_TEXT ADD
#
# The "TOPIC TABLE" block on Sheet2 gains a new column-definition row
# (TOPIC_TEXT / VARCHAR(3000) / NOT NULL) right after TOPIC_NAME and before
# RULE_REFERENCES. Inserting a row shifts the table down by one, so the
# stray blank spacer row that used to separate TOPIC TABLE from RULE_TABLE
# is removed afterwards to keep RULE_TABLE anchored on row 16, exactly as
# it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Insert a new row above the current "RULE_REFERENCES" row (row 12), pushing
# RULE_REFERENCES / CREATED_DATE / the blank spacer row down by one.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the TOPIC_TEXT column definition.
$ws.Range("A12").Value = "TOPIC_TEXT"
$ws.Range("B12").Value = "VARCHAR(3000)"
$ws.Range("C12").Value = "NOT NULL"

# Match the row height used throughout the rest of the table.
$ws.Rows.Item(12).RowHeight = 21

# Remove the now-duplicated blank spacer row (originally row 14, shifted to
# row 15 by the insert above) so RULE_TABLE's header lands back on row 16.
$ws.Rows.Item(15).Delete()
